# Merge the two runs " storing the data any person would " and
# "have (e.g. their name). " (which sit back-to-back inside the
# "Both can inherit from a base class person ..." paragraph) into a
# single run, leaving the rest of the paragraph (and presentation)
# untouched.

$p = $ppt.ActivePresentation

$needle = " storing the data any person would have (e.g. their name). "

$targetShape = $null
$targetParaIdx = -1

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $tf = $shape.TextFrame
            if ($tf.HasText) {
                $tr = $tf.TextRange
                $pCount = $tr.Paragraphs().Count
                for ($pi = 1; $pi -le $pCount; $pi++) {
                    $para = $tr.Paragraphs($pi, 1)
                    if ($para.Text.IndexOf("storing the data any person would") -ge 0) {
                        $targetShape = $shape
                        $targetParaIdx = $pi
                    }
                }
            }
        }
    }
}

if ($targetShape -eq $null) {
    throw "Could not locate the paragraph containing 'storing the data any person would'."
}

$tr = $targetShape.TextFrame.TextRange
$para = $tr.Paragraphs($targetParaIdx, 1)
$offset = $para.Text.IndexOf($needle)
$run = $tr.Characters($para.Start + $offset, $needle.Length)
$run.Text = $needle
